$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Buying Opportunity) updates
$ws.Range("B2").Value  = "NSE:ARIHANTSUP"
$ws.Range("B3").Value  = "NSE:BUTTERFLY"
$ws.Range("B4").Value  = "NSE:INOXWIND"
$ws.Range("B5").Value  = "NSE:IRISDOREME"
$ws.Range("B6").Value  = "NSE:JUBLPHARMA"
$ws.Range("B7").Value  = "NSE:KAMOPAINTS"
$ws.Range("B8").Value  = "NSE:LAOPALA"
$ws.Range("B9").Value  = "NSE:MONTECARLO"
$ws.Range("B10").Value = "NSE:NH"

# These rows no longer have a "Buying Opportunity" entry
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("B19").Value = ""

# Column C (support Zone) updates
$ws.Range("C3").Value  = "NSE:ADFFOODS"
$ws.Range("C4").Value  = "NSE:ARE&M"
$ws.Range("C5").Value  = "NSE:CCL"
$ws.Range("C6").Value  = "NSE:CEATLTD"
$ws.Range("C7").Value  = "NSE:CERA"
$ws.Range("C8").Value  = "NSE:COCHINSHIP"
$ws.Range("C9").Value  = "NSE:GATEWAY"
$ws.Range("C10").Value = "NSE:GEOJITFSL"
$ws.Range("C11").Value = "NSE:GPIL"
$ws.Range("C12").Value = "NSE:HEALTHY"
$ws.Range("C13").Value = "NSE:HIL"
$ws.Range("C14").Value = "NSE:INDORAMA"
$ws.Range("C15").Value = "NSE:KAYNES"
$ws.Range("C16").Value = "NSE:LUXIND"
$ws.Range("C17").Value = "NSE:MASTEK"
$ws.Range("C18").Value = "NSE:NIITLTD"
$ws.Range("C19").Value = "NSE:RAMCOSYS"

# Column E (Short buildup) updates
$ws.Range("E2").Value = "NSE:ADANIENT"
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""

# New row 20
$ws.Range("A20").Value = 18
$ws.Range("C20").Value = "NSE:ROSSELLIND"

# Give the new index cell (A20) the same look as the rest of column A (s="1")
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = 0
